$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: update K4 value (L4 formula is already =100*(K4-$G$4)/K4, recalculates automatically)
$ws.Range("K4").Value = 158.1

# Row 5: fill K5 and fix L5 formula to reference its own row ($G$5 instead of shared $G$4)
$ws.Range("K5").Value = 40.8
$ws.Range("L5").Formula = "=100*(K5-`$G`$5)/K5"

# Row 6: fill K6 and give L6 its own (non-shared) formula
$ws.Range("K6").Value = 38.66
$ws.Range("L6").Formula = "=100*(K6-`$G`$6)/K6"

# Row 7: fill K7 and give L7 its own (non-shared) formula
$ws.Range("K7").Value = 9.152
$ws.Range("L7").Formula = "=100*(K7-`$G`$7)/K7"

# Row 8: fill K8 and give L8 its own (non-shared) formula
$ws.Range("K8").Value = 15.734
$ws.Range("L8").Formula = "=100*(K8-`$G`$8)/K8"

# Row 9: update K9 value (L9 formula already references its own row, recalculates automatically)
$ws.Range("K9").Value = 3.653

# Row 10: fill K10 and give L10 its own (non-shared) formula
$ws.Range("K10").Value = 12.5
$ws.Range("L10").Formula = "=100*(K10-`$G`$10)/K10"

# Update the selected cell in the sheet view
$ws.Range("K6").Select()
